$wb = $excel.ActiveWorkbook

# Sheet "展览" - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 251
$ws1.Range("F9").Value = 269
$ws1.Range("F12").Value = 649
$ws1.Range("F13").Value = 751
$ws1.Range("F14").Value = 1504
$ws1.Range("F15").Value = 1504
$ws1.Range("F17").Value = 28
$ws1.Range("F20").Value = 290
$ws1.Range("F24").Value = 6557
$ws1.Range("F25").Value = 4891
$ws1.Range("F26").Value = 142
$ws1.Range("F29").Value = 155
$ws1.Range("F32").Value = 1275
$ws1.Range("F35").Value = 603
$ws1.Range("F37").Value = 1335
$ws1.Range("F38").Value = 240
$ws1.Range("F40").Value = 143
$ws1.Range("F42").Value = 90
$ws1.Range("F43").Value = 96

# Sheet "本地生活" - column F updates
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2443
$ws3.Range("F5").Value = 49

# Sheet "全部类型" - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 251
$ws4.Range("F8").Value = 49
$ws4.Range("F13").Value = 269
$ws4.Range("F17").Value = 649
$ws4.Range("F18").Value = 751
$ws4.Range("F19").Value = 1504
$ws4.Range("F20").Value = 1504
$ws4.Range("F22").Value = 28
$ws4.Range("F25").Value = 290
$ws4.Range("F30").Value = 6557
$ws4.Range("F31").Value = 4891
$ws4.Range("F32").Value = 142
$ws4.Range("F34").Value = 1275
$ws4.Range("F38").Value = 603
$ws4.Range("F43").Value = 1335
$ws4.Range("F44").Value = 240
$ws4.Range("F45").Value = 143
$ws4.Range("F47").Value = 90
$ws4.Range("F48").Value = 96
